# Set hasChildren (column J) to "Yes" for the rows identified as having children
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(2,6,10,12,18,24,26,29,46,50,53,58,59,60,61,68,70,89,99,114,123,137,142,144,145,150,156,159,166,179,186,191,202,211,217,240,244,266,268,269,270,277,280,282,288,294,302,304,305,324,331,332,342,347,356,363,364,367,374,383,389,395,405,408,412,417,422,429,438,439,444,445,448,449,455,466,469,471,474,476,477,480,485,490,501)

foreach ($r in $rows) {
    $ws.Cells.Item($r, 10).Value = "Yes"
}
